# Implements "1D Lateral interactions implemented, first iteration" edit:
# - Inserts new Mualem Van Genuchten "water" / "clay" / "peat" parameter rows
#   alongside the existing "sand" / "silt" rows (alpha_*, n_*, residual_wc_*).
# - Updates alpha_sand value (4 -> 4.06).
# - Moves the Darcy_friction_factor / tortuosity_air rows down to make room
#   (old rows 46/47 -> new rows 55/56), leaving a blank separator row (54).
# - Adds a handful of blank, number-formatted "scratch" cells in column N
#   (rows 12,13,15,16,17,18,26) mirroring the author's original edit.
# - Updates the active selection to M22 (and drops the old frozen
#   top-left/selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the whole "Mualem Van Genuchten" block (old rows 39-47, which will be
# repopulated at new rows 39-56) so no stale values/styles from the old
# layout (e.g. old C46/C47 text, old B47's number format) bleed through at
# their old positions.
$ws.Range("A39:D47").Clear()

# ---------------------------------------------------------------------
# 1. Mualem Van Genuchten block (rows 39-53 in the final layout)
# ---------------------------------------------------------------------
# Final desired A/B/D values, written directly (the final layout is known
# exactly, so cells are populated top-down rather than simulated via
# row-insert operations).

$ws.Range("A39").Value = "alpha_water"
$ws.Range("B39").Value = 400

$ws.Range("A40").Value = "alpha_sand"
$ws.Range("B40").Value = 4.06

$ws.Range("A41").Value = "alpha_silt"
$ws.Range("B41").Value = 0.65

$ws.Range("A42").Value = "alpha_clay"
$ws.Range("B42").Value = 1.49

$ws.Range("A43").Value = "alpha_peat"
$ws.Range("B43").Value = 2.31
$ws.Range("D43").Value = "from Hydraulic properties of fen peat soils in Poland, Gnatowski 2010"

$ws.Range("A44").Value = "n_water"
$ws.Range("B44").Value = 2.5

$ws.Range("A45").Value = "n_sand"
$ws.Range("B45").Value = 2

$ws.Range("A46").Value = "n_silt"
$ws.Range("B46").Value = 1.7

$ws.Range("A47").Value = "n_clay"
$ws.Range("B47").Value = 1.25

$ws.Range("A48").Value = "n_peat"
$ws.Range("B48").Value = 1.29

$ws.Range("A49").Value = "residual_wc_water"
$ws.Range("B49").Value = 0

$ws.Range("A50").Value = "residual_wc_sand"
$ws.Range("B50").Value = 0

$ws.Range("A51").Value = "residual_wc_silt"
$ws.Range("B51").Value = 0

$ws.Range("A52").Value = "residual_wc_clay"
$ws.Range("B52").Value = 0

$ws.Range("A53").Value = "residual_wc_peat"
$ws.Range("B53").Value = 0

# Row 54 stays empty (separator), matching the original sheet's blank-row
# style (row 45/54 before, row 38/54 pattern).

# ---------------------------------------------------------------------
# 2. Move Darcy_friction_factor / tortuosity_air down to rows 55-56
# ---------------------------------------------------------------------

$ws.Range("A55").Value = "Darcy_friction_factor"
$ws.Range("B55").Value = 0.1
$ws.Range("C55").Value = "rough-pipe regime"

$ws.Range("A56").Value = "tortuosity_air"
$ws.Range("B56").Value = 2.5
$ws.Range("B56").NumberFormat = "General"
$ws.Range("C56").Value = "used in Carman Kozeny model"

# Clear out the old rows 46/47 contents that are no longer at their old
# position (the block above already overwrote A46:D47 with the new
# water/clay data, so nothing further to clear there).

# ---------------------------------------------------------------------
# 3. Small formatted "scratch" cells added in column N
# ---------------------------------------------------------------------

$ws.Range("N12").NumberFormat = "0.00E+00"
$ws.Range("N13").NumberFormat = "0.00E+00"
$ws.Range("N15").NumberFormat = "0.00E+00"
$ws.Range("N16").NumberFormat = "0.00E+00"
$ws.Range("N17").NumberFormat = "0.00E+00"
$ws.Range("N18").NumberFormat = "0.00E+00"
$ws.Range("N26").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------
# 4. Selection / view state
# ---------------------------------------------------------------------

$ws.Range("M22").Select()
